$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("B1").Value = "Train Loss"
$ws.Range("C1").Value = "Test Loss"
$ws.Range("D1").Value = "Accuracy"

# Match the existing header formatting (bold font, thin border, centered
# horizontal / top vertical alignment) on the two newly introduced header
# cells by copying the format already used for A1/"Epoch".
$ws.Range("A1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New train-loss / test-loss / accuracy values for epochs 1-10.
$data = @(
  @(1.6553, 1.4902, 46.0594),
  @(1.4175, 1.3641, 49.9046),
  @(1.3171, 1.2835, 52.0862),
  @(1.2559, 1.2457, 53.1497),
  @(1.2144, 1.2281, 53.9951),
  @(1.1906, 1.2070, 53.8042),
  @(1.1697, 1.1787, 56.2040),
  @(1.1571, 1.1751, 55.0586),
  @(1.1491, 1.1693, 55.6586),
  @(1.1432, 1.1663, 55.3313)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $vals = $data[$i]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Remove the old rows for epochs 11-20 (previously rows 12-21), which no
# longer exist in the trimmed table.
$ws.Range("A12:D21").EntireRow.Delete()
